$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:H2")
$rng.Font.Name = "Calibri"
$rng.Font.Size = 10

$ws.Range("A2").Value = "MCH216"
$ws.Range("C2").Value = "VARIOUS LEAFLETS, PAMPHLETS, PRESS CUTTINGS, LETTERS ETC, INCLUDING MATERIAL ON THE DUTCH BOA- URBAN CAMPAIGN "
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1C | GRAP COUNT NUMER: NONE"
$ws.Range("H2").Value = ""
